$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 955.3333
$ws.Range("I15").Value = 955.3333
$ws.Range("K15").Value = 2865.9999
$ws.Range("M15").Value = -2696.9999
$ws.Range("H41").Value = 1105.4117
$ws.Range("J41").Value = 537.3333
$ws.Range("L41").Value = 537.3333
$ws.Range("N41").Value = -1417.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 784.1429000000001
$ws.Range("I2").Value = 707.5789
$ws.Range("J2").Value = 945.7778
$ws.Range("K2").Value = 707.5789
$ws.Range("L2").Value = 945.7778
$ws.Range("M2").Value = -594.5789
$ws.Range("N2").Value = -1171.7778
$ws.Range("H32").Value = 12512.103
$ws.Range("I32").Value = 12921.257
$ws.Range("K32").Value = 12921.257
$ws.Range("M32").Value = -12634.257
$ws.Range("H45").Value = 2005
$ws.Range("I45").Value = 1981.8182
$ws.Range("J45").Value = 2056
$ws.Range("K45").Value = 1981.8182
$ws.Range("L45").Value = 2056
$ws.Range("M45").Value = -1604.8182
$ws.Range("N45").Value = -2810
$ws.Range("H61").Value = 13891711
$ws.Range("I61").Value = 17859742
$ws.Range("K61").Value = 17859742
$ws.Range("M61").Value = -17859530
$ws.Range("H97").Value = 7857.7144
$ws.Range("I97").Value = 8369
$ws.Range("J97").Value = 1211
$ws.Range("K97").Value = 8369
$ws.Range("L97").Value = 1211
$ws.Range("M97").Value = -7873
$ws.Range("N97").Value = -2203
$ws.Range("H106").Value = 50696.668
$ws.Range("J106").Value = 50696.668
$ws.Range("L106").Value = 50696.668
$ws.Range("N106").Value = -53220.668
$ws.Range("H110").Value = 1631.0526
$ws.Range("I110").Value = 1493.125
$ws.Range("K110").Value = 1493.125
$ws.Range("M110").Value = 551.875
$ws.Range("H116").Value = 784.1429000000001
$ws.Range("I116").Value = 707.5789
$ws.Range("J116").Value = 945.7778
$ws.Range("K116").Value = 707.5789
$ws.Range("L116").Value = 945.7778
$ws.Range("M116").Value = 1586.4211
$ws.Range("N116").Value = -5533.7778
$ws.Range("H122").Value = 6505.125
$ws.Range("I122").Value = 7753.8887
$ws.Range("J122").Value = 2758.8333
$ws.Range("K122").Value = 23261.6661
$ws.Range("L122").Value = 8276.499899999999
$ws.Range("M122").Value = -20811.6661
$ws.Range("N122").Value = -13176.4999
$ws.Range("H132").Value = 7144823.5
$ws.Range("I132").Value = 10871085
$ws.Range("J132").Value = 2821
$ws.Range("K132").Value = 32613255
$ws.Range("L132").Value = 8463
$ws.Range("M132").Value = -32610725
$ws.Range("N132").Value = -13523
$ws.Range("H136").Value = 13891711
$ws.Range("I136").Value = 17859742
$ws.Range("K136").Value = 53579226
$ws.Range("M136").Value = -53576676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 784.1429000000001
$ws.Range("I3").Value = 707.5789
$ws.Range("J3").Value = 945.7778
$ws.Range("K3").Value = 707.5789
$ws.Range("L3").Value = 945.7778
$ws.Range("M3").Value = -593.5789
$ws.Range("N3").Value = -1173.7778
$ws.Range("H11").Value = 980.2727
$ws.Range("I11").Value = 964
$ws.Range("J11").Value = 999.8
$ws.Range("K11").Value = 964
$ws.Range("L11").Value = 999.8
$ws.Range("M11").Value = -824
$ws.Range("N11").Value = -1279.8
$ws.Range("H16").Value = 5004
$ws.Range("I16").Value = 5004
$ws.Range("K16").Value = 5004
$ws.Range("M16").Value = -4834
$ws.Range("H94").Value = 1164.4762
$ws.Range("I94").Value = 1021.1539
$ws.Range("J94").Value = 1397.375
$ws.Range("K94").Value = 1021.1539
$ws.Range("L94").Value = 1397.375
$ws.Range("M94").Value = -570.1539
$ws.Range("N94").Value = -2299.375
$ws.Range("H105").Value = 3824.468
$ws.Range("I105").Value = 2636.95
$ws.Range("J105").Value = 4704.1113
$ws.Range("K105").Value = 2636.95
$ws.Range("L105").Value = 4704.1113
$ws.Range("M105").Value = -889.9499999999998
$ws.Range("N105").Value = -8198.1113
$ws.Range("H107").Value = 2120.5833
$ws.Range("I107").Value = 2168.9473
$ws.Range("J107").Value = 1936.8
$ws.Range("K107").Value = 2168.9473
$ws.Range("L107").Value = 1936.8
$ws.Range("M107").Value = -248.9472999999998
$ws.Range("N107").Value = -5776.8
$ws.Range("H134").Value = 3636.7334
$ws.Range("I134").Value = 2253.3333
$ws.Range("J134").Value = 5711.8335
$ws.Range("K134").Value = 6759.999899999999
$ws.Range("L134").Value = 17135.5005
$ws.Range("M134").Value = -4224.999899999999
$ws.Range("N134").Value = -22205.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 941.16
$ws.Range("I16").Value = 795.4666999999999
$ws.Range("J16").Value = 1159.7
$ws.Range("K16").Value = 795.4666999999999
$ws.Range("L16").Value = 1159.7
$ws.Range("M16").Value = -508.4666999999999
$ws.Range("N16").Value = -1733.7
$ws.Range("H19").Value = 123.85714
$ws.Range("I19").Value = 131
$ws.Range("J19").Value = 106
$ws.Range("K19").Value = 131
$ws.Range("L19").Value = 106
$ws.Range("M19").Value = 39
$ws.Range("N19").Value = -446
$ws.Range("H24").Value = 123.85714
$ws.Range("I24").Value = 131
$ws.Range("J24").Value = 106
$ws.Range("K24").Value = 131
$ws.Range("L24").Value = 106
$ws.Range("M24").Value = 39
$ws.Range("N24").Value = -446
$ws.Range("H113").Value = 941.16
$ws.Range("I113").Value = 795.4666999999999
$ws.Range("J113").Value = 1159.7
$ws.Range("K113").Value = 795.4666999999999
$ws.Range("L113").Value = 1159.7
$ws.Range("M113").Value = 1374.5333
$ws.Range("N113").Value = -5499.7
$ws.Range("H122").Value = 1552.3334
$ws.Range("I122").Value = 1679.8667
$ws.Range("J122").Value = 1392.9166
$ws.Range("K122").Value = 5039.6001
$ws.Range("L122").Value = 4178.7498
$ws.Range("M122").Value = -2589.6001
$ws.Range("N122").Value = -9078.7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 578.2917
$ws.Range("I4").Value = 174
$ws.Range("J4").Value = 2599.75
$ws.Range("K4").Value = 522
$ws.Range("L4").Value = 7799.25
$ws.Range("M4").Value = -410
$ws.Range("N4").Value = -8023.25
$ws.Range("H131").Value = 903.125
$ws.Range("I131").Value = 825.9
$ws.Range("J131").Value = 913.02563
$ws.Range("K131").Value = 2477.7
$ws.Range("L131").Value = 2739.07689
$ws.Range("M131").Value = 2562.3
$ws.Range("N131").Value = -12819.07689

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 1929.6666
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 2494.5
$ws.Range("K25").Value = 800
$ws.Range("L25").Value = 2494.5
$ws.Range("M25").Value = -271
$ws.Range("N25").Value = -3552.5
$ws.Range("H97").Value = 1356
$ws.Range("I97").Value = 1483.875
$ws.Range("J97").Value = 1100.25
$ws.Range("K97").Value = 1483.875
$ws.Range("L97").Value = 1100.25
$ws.Range("M97").Value = -987.875
$ws.Range("N97").Value = -2092.25
$ws.Range("H102").Value = 2186.5833
$ws.Range("I102").Value = 2561.6538
$ws.Range("K102").Value = 2561.6538
$ws.Range("M102").Value = -939.6538
$ws.Range("H126").Value = 4018.8
$ws.Range("I126").Value = 2536.3076
$ws.Range("J126").Value = 5624.8335
$ws.Range("K126").Value = 7608.9228
$ws.Range("L126").Value = 16874.5005
$ws.Range("M126").Value = -5138.9228
$ws.Range("N126").Value = -21814.5005
$ws.Range("H132").Value = 4692.3213
$ws.Range("I132").Value = 3491.5
$ws.Range("J132").Value = 5893.143
$ws.Range("K132").Value = 10474.5
$ws.Range("L132").Value = 17679.429
$ws.Range("M132").Value = -7944.5
$ws.Range("N132").Value = -22739.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5242.727
$ws.Range("J40").Value = 5537.375
$ws.Range("L40").Value = 5537.375
$ws.Range("N40").Value = -5809.375
$ws.Range("H61").Value = 1247.1875
$ws.Range("I61").Value = 1139.6428
$ws.Range("K61").Value = 1139.6428
$ws.Range("M61").Value = -937.6428000000001
$ws.Range("H68").Value = 1963.5294
$ws.Range("I68").Value = 1741.6666
$ws.Range("J68").Value = 2496
$ws.Range("K68").Value = 1741.6666
$ws.Range("L68").Value = 2496
$ws.Range("M68").Value = -992.6666
$ws.Range("N68").Value = -3994
$ws.Range("H71").Value = 1963.5294
$ws.Range("I71").Value = 1741.6666
$ws.Range("J71").Value = 2496
$ws.Range("K71").Value = 8708.333000000001
$ws.Range("L71").Value = 12480
$ws.Range("M71").Value = -4964.333000000001
$ws.Range("N71").Value = -19968
$ws.Range("H81").Value = 37759.223
$ws.Range("J81").Value = 37759.223
$ws.Range("L81").Value = 37759.223
$ws.Range("N81").Value = -39755.223
$ws.Range("H84").Value = 37759.223
$ws.Range("J84").Value = 37759.223
$ws.Range("L84").Value = 113277.669
$ws.Range("N84").Value = -123261.669
$ws.Range("H93").Value = 1590.8
$ws.Range("I93").Value = 1700
$ws.Range("J93").Value = 1551.091
$ws.Range("K93").Value = 1700
$ws.Range("L93").Value = 1551.091
$ws.Range("M93").Value = -452
$ws.Range("N93").Value = -4047.091
$ws.Range("H113").Value = 1247.1875
$ws.Range("I113").Value = 1139.6428
$ws.Range("K113").Value = 1139.6428
$ws.Range("M113").Value = 1030.3572
$ws.Range("H122").Value = 4971.6665
$ws.Range("I122").Value = 5270.8335
$ws.Range("J122").Value = 4672.5
$ws.Range("K122").Value = 15812.5005
$ws.Range("L122").Value = 14017.5
$ws.Range("M122").Value = -13362.5005
$ws.Range("N122").Value = -18917.5
$ws.Range("H132").Value = 8339796
$ws.Range("I132").Value = 4135.684
$ws.Range("J132").Value = 22737754
$ws.Range("K132").Value = 12407.052
$ws.Range("L132").Value = 68213262
$ws.Range("M132").Value = -9877.052
$ws.Range("N132").Value = -68218322
$ws.Range("H136").Value = 71450270
$ws.Range("I136").Value = 125003470
$ws.Range("J136").Value = 46001.668
$ws.Range("K136").Value = 375010410
$ws.Range("L136").Value = 138005.004
$ws.Range("M136").Value = -375007860
$ws.Range("N136").Value = -143105.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1573.95
$ws.Range("I113").Value = 532.4167
$ws.Range("J113").Value = 3136.25
$ws.Range("K113").Value = 1597.2501
$ws.Range("L113").Value = 9408.75
$ws.Range("M113").Value = 572.7499
$ws.Range("N113").Value = -13748.75
$ws.Range("H122").Value = 2757.95
$ws.Range("I122").Value = 2714.1538
$ws.Range("J122").Value = 2839.2856
$ws.Range("K122").Value = 8142.4614
$ws.Range("L122").Value = 8517.856800000001
$ws.Range("M122").Value = -5692.4614
$ws.Range("N122").Value = -13417.8568
$ws.Range("H132").Value = 2200.5
$ws.Range("I132").Value = 1167
$ws.Range("K132").Value = 3501
$ws.Range("M132").Value = -971
$ws.Range("H136").Value = 1572.7368
$ws.Range("I136").Value = 1678.7142
$ws.Range("J136").Value = 1276
$ws.Range("K136").Value = 5036.142599999999
$ws.Range("L136").Value = 3828
$ws.Range("M136").Value = -2486.142599999999
$ws.Range("N136").Value = -8928
